# ---------------------------------------------------------------------------
# PlayerPerformance_5957.xlsx update
#   1. Insert a brand-new "Player Info" sheet in front of the existing
#      "ODI Batting" / "ODI Bowling" sheets, with the player's basic info.
#   2. On "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE and replace the
#      full scorecard URLs with just the trailing MatchCode number.
#   3. On "ODI Bowling": same MATCH_CARD_LINK -> MATCH_CODE rename + URL ->
#      MatchCode replacement.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before "ODI Batting" (becomes sheet 1)
#    NOTE: worksheet handles obtained before a sheet Add() can end up bound
#    to the wrong position afterwards, so we re-fetch sheets by name AFTER
#    any structural change (insert/delete/move) instead of reusing old refs.
# ---------------------------------------------------------------------------
$odiBattingBeforeInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($odiBattingBeforeInsert)
$playerInfo.Name = "Player Info"

# Match the page-setup / outline conventions already used by the rest of
# the workbook (new sheets otherwise get iron's own defaults).
$playerInfo.Outline.SummaryRow = 1
$playerInfo.Outline.SummaryColumn = 1
$playerInfo.PageSetup.LeftMargin = 54
$playerInfo.PageSetup.RightMargin = 54
$playerInfo.PageSetup.TopMargin = 72
$playerInfo.PageSetup.BottomMargin = 72
$playerInfo.PageSetup.HeaderMargin = 36
$playerInfo.PageSetup.FooterMargin = 36

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$piHeader = $playerInfo.Range("A1:D1")
$piHeader.Font.Bold = $true
$piHeader.Borders.LineStyle = 1
$piHeader.HorizontalAlignment = -4108
$piHeader.VerticalAlignment = -4160

# ID is stored as text (e.g. "5957"), not a number, so force text format
# before writing the value, then drop back to the default style so the
# cell isn't left with an unwanted "text number format" style.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "5957"
$playerInfo.Range("A2").Style = "Normal"

$playerInfo.Range("B2").Value = "Nasum Ahmed"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# ---------------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code
#    (re-fetch the sheet now that the workbook layout changed above)
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"

$battingRows  = @(2, 3, 4, 5, 6, 7, 8)
$battingCodes = @("4606", "4611", "4616", "4682", "4726", "4729", "4734")
for ($i = 0; $i -lt $battingRows.Length; $i++) {
    $cell = $odiBatting.Range("D" + $battingRows[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$i]
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code
# ---------------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Range("B1").Value = "MATCH_CODE"

$bowlingRows  = @(2, 3, 4, 5, 6, 7)
$bowlingCodes = @("4606", "4611", "4616", "4682", "4726", "4734")
for ($i = 0; $i -lt $bowlingRows.Length; $i++) {
    $cell = $odiBowling.Range("B" + $bowlingRows[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$i]
    $cell.Style = "Normal"
}
